$wb = $excel.ActiveWorkbook
$croatia = $wb.Worksheets.Item("Croatia")

# Duplicate Croatia's sheet, placing the copy right after Croatia.
$croatia.Copy($null, $croatia)

$greece = $wb.Worksheets.Item($croatia.Index + 1)
$greece.Name = "Greece"

# New Jira ticket strings for this commit.
$greece.Range("B4").Value = "NGC-4119/T3196"
$croatia.Range("B4").Value = "NGC-3139/T2073"

# Selections / active-cell bookkeeping that Excel records per sheet.
$greece.Range("B4").Select()
$croatia.Range("E14").Select()

$greece.Activate()
